$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.803.17'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.740.25'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '405.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.32'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.722.17'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.718'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.165'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -10.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000364'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -7.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.32'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.284.97'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.60'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.41'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +11.80%  '
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.724.34'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '65.906.37'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('E21').Value = '  -6.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.72'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -9.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.40'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -7.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.01'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '36.14'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.52'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +10.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.07'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -10.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.39'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('E32').Value = '  -3.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.04'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.05%  '
$ws.Range('E34').Value = '  -4.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '38.41'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -8.15%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.96'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0735'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0450'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -8.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.996'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.81'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.24%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.19'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +21.79%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.134'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -8.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '145.03'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.21'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.81%  '
$ws.Range('E47').Value = '  -4.08%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.22'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.78'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.54'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.286'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.66%  '
